# Fix Training Data Issue (#48)
# The "Date" column (BF) on Sheet1 contains the literal text "5-20-2012-13"
# for every data row (rows 2-31). That string is wrong / ambiguous (it is
# "<month>-<day>-<season>"), so replace it with the correct ISO-ish date
# string "2013-05-20" for the 5/20/2013 game date, keeping the value as
# plain text (not an Excel date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "5-20-2012-13"
$newValue = "2013-05-20"

$firstRow = 2
$lastRow = 31
$col = 58  # column BF

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)

    if ($cell.Value2 -ne $oldValue) {
        continue
    }

    # Writing a plain "yyyy-mm-dd"-shaped string straight into .Value2 makes
    # Excel "helpfully" reinterpret it as a date serial number. Force the
    # cell to Text first so the string round-trips unchanged, then restore
    # the cell's style (so no lingering number-format is left behind) while
    # keeping the literal text value.
    $cell.NumberFormat = "@"
    $cell.Value2 = $newValue
    $cell.Style = "Normal"
}
